$d = $word.ActiveDocument

# Find the paragraph containing the first bullet item and grab its paragraph object
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Регистрация и авторизация пользователей") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Insert a new paragraph right after it, matching the same list style
$newRange = $target.Range.InsertParagraphAfter()

# The newly created paragraph is the next paragraph after $target
$newPara = $target.Next()
$newPara.Style = $target.Style
$newPara.Range.ListFormat.ApplyBulletDefault()

$newPara.Range.Text = "Изменение данных профиля"
$newPara.Range.Font.Name = "Times New Roman"
$newPara.Range.Font.Size = 16
